$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# Add a new skill entry in row 14, mirroring the formatting used by the
# preceding data rows, reflecting BUFF and Receive Damage context
# triggers. Copy the row formatting from the row above (row 13) so the
# new row matches the existing "data row" look exactly.
$ws.Range("A13:E13").Copy()
$ws.Range("A14:E14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 2).Value = "ScaleOneself"
$ws.Cells.Item(14, 3).Value = "StatusSkill"
$ws.Cells.Item(14, 4).Value = 10
$ws.Cells.Item(14, 5).Value = 3

# Update the active selection to reflect where the user ended up after
# the edit.
$ws.Range("F17").Select()
